$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.827023400069947
$ws.Range("C2").Value = 0.4888028508193502
$ws.Range("D2").Value = 1.471596344594558
$ws.Range("B3").Value = 0.06190938385428302
$ws.Range("C3").Value = 2.26017545783639
$ws.Range("B4").Value = 2.750952690361402
$ws.Range("B5").Value = 0.7320287013142677
$ws.Range("C5").Value = 1.522302341005567
$ws.Range("D5").Value = 1.467338125109618
$ws.Range("E5").Value = 1.793638825175652
$ws.Range("B6").Value = 1.9795770604465
$ws.Range("C6").Value = 1.392536339603099
$ws.Range("D6").Value = 1.481510914913067
$ws.Range("B7").Value = 1.30090756340108
$ws.Range("C7").Value = 1.963500078556037
$ws.Range("B8").Value = 2.012737994557474
$ws.Range("B9").Value = 0.8535918672211444
$ws.Range("C9").Value = 1.782613822028589
$ws.Range("D9").Value = 0.7889205787030562
$ws.Range("E9").Value = 0.6461210271256811
$ws.Range("B10").Value = 1.972402569862437
$ws.Range("C10").Value = 0.8237969262448452
$ws.Range("D10").Value = 0.4993394432980468
$ws.Range("B11").Value = 0.796001625707733
$ws.Range("C11").Value = 0.9151352642427935
$ws.Range("B12").Value = 1.082826669985442
$ws.Range("B13").Value = 0.573129683809814
$ws.Range("C13").Value = 0.1945252858301101
$ws.Range("D13").Value = 0.2603721808367071
$ws.Range("B14").Value = 0.1824667929082922
$ws.Range("C14").Value = 0.4481552108942597
$ws.Range("B15").Value = 0.3345343741504182
